# Regenerate save_data column G ("K") values for the 2023 lynn_lance sheet.
# The commit regenerated K (strikeouts) from the source log instead of the
# old "Strike#" figure, so only the G2:G38 cell values need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 4
    4  = 7
    5  = 6
    6  = 3
    7  = 1
    8  = 1
    9  = 1
    10 = 3
    11 = 9
    12 = 9
    13 = 7
    14 = 5
    15 = 7
    16 = 6
    17 = 11
    18 = 7
    19 = 8
    20 = 16
    21 = 6
    22 = 4
    23 = 3
    24 = 5
    25 = 6
    26 = 7
    27 = 4
    28 = 8
    29 = 10
    30 = 4
    31 = 8
    32 = 11
    33 = 6
    34 = 6
    35 = 6
    36 = 5
    37 = 3
    38 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
